$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(10).Insert()

$ws.Range("A10").Value = 3
$ws.Range("B10").Value = "Femacal de La Calera"
$ws.Range("C10").Value = "Coquimbo"
$ws.Range("D10").Value = 44819
$ws.Range("D10").NumberFormat = $ws.Range("D11").NumberFormat
$ws.Range("E10").Value = 5
$ws.Range("F10").Value = 100112035
$ws.Range("G10").Value = "Bruselas (repollito)"
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 45
$ws.Range("K10").Value = 16000
$ws.Range("L10").Value = 16000
$ws.Range("M10").Value = 16000
$ws.Range("N10").Value = "`$/malla 15 kilos"
$ws.Range("O10").Value = "Provincia de Quillota"
$ws.Range("P10").Value = 1067
$ws.Range("Q10").Value = 15
$ws.Range("R10").Value = "Hortaliza"

Write-Output "done"
